$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly Jengibre price table gained a new entry at the top (row 9),
# pushing every existing data row down by one. What used to be row 50
# now becomes row 51, so first extend the sheet with a new row 51 that
# duplicates the constant columns from row 50 (everything except the
# Fecha/Volumen/Precio columns, which differ per row).
$ws.Range("A51").Value = $ws.Range("A50").Value2
$ws.Range("B51").Value = $ws.Range("B50").Value2
$ws.Range("C51").Value = $ws.Range("C50").Value2
$ws.Range("E51").Value = $ws.Range("E50").Value2
$ws.Range("F51").Value = $ws.Range("F50").Value2
$ws.Range("G51").Value = $ws.Range("G50").Value2
$ws.Range("H51").Value = $ws.Range("H50").Value2
$ws.Range("I51").Value = $ws.Range("I50").Value2
$ws.Range("N51").Value = $ws.Range("N50").Value2
$ws.Range("O51").Value = $ws.Range("O50").Value2
$ws.Range("Q51").Value = $ws.Range("Q50").Value2
$ws.Range("R51").Value = $ws.Range("R50").Value2

# Match the date-time number format used by the rest of column D.
$ws.Range("D51").NumberFormat = $ws.Range("D50").NumberFormat

# Now update Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), and Precio $/Kg (P) for rows 9-51:
# row 9 receives a brand-new weekly entry, rows 10-51 take over the values
# that used to belong to rows 9-50 (the table is shifted down by one row).

$ws.Range("D9").Value = 44487
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14150
$ws.Range("P9").Value = 1088

$ws.Range("D10").Value = 44340
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 23000
$ws.Range("M10").Value = 21275
$ws.Range("P10").Value = 1637

$ws.Range("D11").Value = 44376
$ws.Range("J11").Value = 580
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13103
$ws.Range("P11").Value = 1008

$ws.Range("D12").Value = 44417
$ws.Range("J12").Value = 230
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 13565
$ws.Range("P12").Value = 1043

$ws.Range("D13").Value = 44245
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 38000
$ws.Range("L13").Value = 40000
$ws.Range("M13").Value = 38850
$ws.Range("P13").Value = 2988

$ws.Range("D14").Value = 44382
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13417
$ws.Range("P14").Value = 1032

$ws.Range("D15").Value = 44284
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 24575
$ws.Range("P15").Value = 1890

$ws.Range("D16").Value = 44315
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 26000
$ws.Range("M16").Value = 25425
$ws.Range("P16").Value = 1956

$ws.Range("D17").Value = 44343
$ws.Range("J17").Value = 290
$ws.Range("K17").Value = 23000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 23897
$ws.Range("P17").Value = 1838

$ws.Range("D18").Value = 44294
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24150
$ws.Range("P18").Value = 1858

$ws.Range("D19").Value = 44356
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15425
$ws.Range("P19").Value = 1187

$ws.Range("D20").Value = 44410
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14575
$ws.Range("P20").Value = 1121

$ws.Range("D21").Value = 44319
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 23000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 24150
$ws.Range("P21").Value = 1858

$ws.Range("D22").Value = 44473
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 12400
$ws.Range("P22").Value = 954

$ws.Range("D23").Value = 44236
$ws.Range("J23").Value = 210
$ws.Range("K23").Value = 47000
$ws.Range("L23").Value = 47000
$ws.Range("M23").Value = 47000
$ws.Range("P23").Value = 3615

$ws.Range("D24").Value = 44168
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 32000
$ws.Range("M24").Value = 31080
$ws.Range("P24").Value = 2391

$ws.Range("D25").Value = 44161
$ws.Range("J25").Value = 330
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29394
$ws.Range("P25").Value = 2261

$ws.Range("D26").Value = 44280
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 23000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 24150
$ws.Range("P26").Value = 1858

$ws.Range("D27").Value = 44445
$ws.Range("J27").Value = 220
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 13000
$ws.Range("M27").Value = 12455
$ws.Range("P27").Value = 958

$ws.Range("D28").Value = 44335
$ws.Range("J28").Value = 170
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 27000
$ws.Range("M28").Value = 25824
$ws.Range("P28").Value = 1986

$ws.Range("D29").Value = 44252
$ws.Range("J29").Value = 130
$ws.Range("K29").Value = 33000
$ws.Range("L29").Value = 35000
$ws.Range("M29").Value = 34077
$ws.Range("P29").Value = 2621

$ws.Range("D30").Value = 44349
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 23000
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = 23800
$ws.Range("P30").Value = 1831

$ws.Range("D31").Value = 44350
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 23000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 24150
$ws.Range("P31").Value = 1858

$ws.Range("D32").Value = 44385
$ws.Range("J32").Value = 220
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13091
$ws.Range("P32").Value = 1007

$ws.Range("D33").Value = 44452
$ws.Range("J33").Value = 290
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 12414
$ws.Range("P33").Value = 955

$ws.Range("D34").Value = 44435
$ws.Range("J34").Value = 580
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = 12500
$ws.Range("P34").Value = 962

$ws.Range("D35").Value = 44242
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 40000
$ws.Range("L35").Value = 42000
$ws.Range("M35").Value = 41200
$ws.Range("P35").Value = 3169

$ws.Range("D36").Value = 44433
$ws.Range("J36").Value = 320
$ws.Range("K36").Value = 12000
$ws.Range("L36").Value = 13000
$ws.Range("M36").Value = 12531
$ws.Range("P36").Value = 964

$ws.Range("D37").Value = 44307
$ws.Range("J37").Value = 580
$ws.Range("K37").Value = 23000
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = 23897
$ws.Range("P37").Value = 1838

$ws.Range("D38").Value = 44263
$ws.Range("J38").Value = 70
$ws.Range("K38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = 30000
$ws.Range("P38").Value = 2308

$ws.Range("D39").Value = 44306
$ws.Range("J39").Value = 230
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = 24435
$ws.Range("P39").Value = 1880

$ws.Range("D40").Value = 44369
$ws.Range("J40").Value = 290
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 17000
$ws.Range("M40").Value = 16172
$ws.Range("P40").Value = 1244

$ws.Range("D41").Value = 44172
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 27000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = 28800
$ws.Range("P41").Value = 2215

$ws.Range("D42").Value = 44301
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 23000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 24200
$ws.Range("P42").Value = 1862

$ws.Range("D43").Value = 44328
$ws.Range("J43").Value = 290
$ws.Range("K43").Value = 23000
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = 23828
$ws.Range("P43").Value = 1833

$ws.Range("D44").Value = 44333
$ws.Range("J44").Value = 80
$ws.Range("K44").Value = 23000
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 24250
$ws.Range("P44").Value = 1865

$ws.Range("D45").Value = 44466
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 13000
$ws.Range("M45").Value = 12400
$ws.Range("P45").Value = 954

$ws.Range("D46").Value = 44270
$ws.Range("J46").Value = 220
$ws.Range("K46").Value = 33000
$ws.Range("L46").Value = 35000
$ws.Range("M46").Value = 34091
$ws.Range("P46").Value = 2622

$ws.Range("D47").Value = 44438
$ws.Range("J47").Value = 300
$ws.Range("K47").Value = 13000
$ws.Range("L47").Value = 14000
$ws.Range("M47").Value = 13333
$ws.Range("P47").Value = 1026

$ws.Range("D48").Value = 44461
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 12000
$ws.Range("L48").Value = 13000
$ws.Range("M48").Value = 12400
$ws.Range("P48").Value = 954

$ws.Range("D49").Value = 44389
$ws.Range("J49").Value = 230
$ws.Range("K49").Value = 13000
$ws.Range("L49").Value = 14000
$ws.Range("M49").Value = 13609
$ws.Range("P49").Value = 1047

$ws.Range("D50").Value = 44312
$ws.Range("J50").Value = 190
$ws.Range("K50").Value = 24000
$ws.Range("L50").Value = 25000
$ws.Range("M50").Value = 24632
$ws.Range("P50").Value = 1895

$ws.Range("D51").Value = 44326
$ws.Range("J51").Value = 290
$ws.Range("K51").Value = 21000
$ws.Range("L51").Value = 23000
$ws.Range("M51").Value = 22172
$ws.Range("P51").Value = 1706
